$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "34.150.92"
$ws.Range("E2").Value = "  -1.56%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.796.56"
$ws.Range("E3").Value = "  +0.48%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.10%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'223.07"
$ws.Range("E5").Value = "  +0.03%  "

# Row 6 - XRP
$ws.Range("D6").Value = "'0.551"
$ws.Range("E6").Value = "  -0.63%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.08%  "

# Row 8 - Solana
$ws.Range("D8").Value = "'32.36"
$ws.Range("E8").Value = "  -0.61%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  +2.25%  "

# Row 10 - Dogecoin
$ws.Range("D10").Value = "'0.0720"
$ws.Range("E10").Value = "  +5.17%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  -1.27%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "2.054.70"
$ws.Range("E12").Value = "  +0.45%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.793.87"
$ws.Range("E13").Value = "  +0.39%  "

# Row 14 - Chainlink
$ws.Range("D14").Value = "'10.75"
$ws.Range("E14").Value = "  -2.18%  "

# Row 15 - Polygon
$ws.Range("E15").Value = "  +0.47%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "34.150.93"
$ws.Range("E16").Value = "  -1.58%  "

# Row 17 - Polkadot
$ws.Range("D17").Value = "'4.22"
$ws.Range("E17").Value = "  -1.55%  "

# Row 18 - Litecoin
$ws.Range("D18").Value = "'68.18"
$ws.Range("E18").Value = "  -0.42%  "

# Row 19 - BitcoinCash
$ws.Range("D19").Value = "'246.69"
$ws.Range("E19").Value = "  -2.47%  "

# Row 20 - ShibaInu
$ws.Range("D20").Value = "0.0₃0788"
$ws.Range("E20").Value = "  +0.60%  "

# Row 21 - Dai
$ws.Range("E21").Value = "  -0.09%  "

# Row 23 - Uniswap
$ws.Range("E23").Value = "  -1.60%  "

# Row 24 - Toncoin
$ws.Range("D24").Value = "'2.12"
$ws.Range("E24").Value = "  -0.54%  "

# Row 25 - Monero
$ws.Range("D25").Value = "'158.94"
$ws.Range("E25").Value = "  +0.11%  "

# Row 26 - EthereumClassic
$ws.Range("D26").Value = "'16.55"
$ws.Range("E26").Value = "  +1.28%  "

# Row 27 - Cosmos
$ws.Range("D27").Value = "'7.07"
$ws.Range("E27").Value = "  +0.45%  "

# Row 29 - BinanceUSD
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.12%  "

# Row 30 - Hedera
$ws.Range("E30").Value = "  +1.34%  "

# Row 31 - was Filecoin, now PancakeSwap
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "'1.21"
$ws.Range("E31").Value = "  +1.95%  "

# Row 32 - was PancakeSwap, now Filecoin
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "'3.72"
$ws.Range("E32").Value = "  -0.61%  "

# Row 33 - InternetComputer(DFINITY)
$ws.Range("D33").Value = "'3.52"
$ws.Range("E33").Value = "  -1.26%  "

# Row 34 - LidoDAOToken
$ws.Range("E34").Value = "  -0.52%  "

# Row 35 - Maker
$ws.Range("D35").Value = "1.414.07"
$ws.Range("E35").Value = "  -1.10%  "

# Row 36 - ImmutableX
$ws.Range("E36").Value = "  +2.15%  "

# Row 37 - TrustWalletToken
$ws.Range("E37").Value = "  +0.08%  "

# Row 38 - VeChain
$ws.Range("D38").Value = "'0.0187"
$ws.Range("E38").Value = "  -1.26%  "

# Row 39 - ARBITRUM
$ws.Range("E39").Value = "  +5.06%  "

# Row 40 - Aave
$ws.Range("D40").Value = "'80.43"
$ws.Range("E40").Value = "  -2.83%  "

# Row 41 - MXToken
$ws.Range("E41").Value = "  -2.79%  "

# Row 42 - HuobiToken
$ws.Range("E42").Value = "  -0.46%  "

# Row 43 - RenderToken
$ws.Range("E43").Value = "  +4.82%  "

# Row 44 - FraxShare
$ws.Range("D44").Value = "'5.97"
$ws.Range("E44").Value = "  +0.51%  "

# Row 45 - Kaspa
$ws.Range("D45").Value = "'0.0497"
$ws.Range("E45").Value = "  -1.40%  "

# Row 46 - Quant
$ws.Range("D46").Value = "'106.99"
$ws.Range("E46").Value = "  +2.94%  "

# Row 47 - RocketPoolETH
$ws.Range("D47").Value = "1.954.80"
$ws.Range("E47").Value = "  +0.71%  "

# Row 48 - WEMIXToken
$ws.Range("E48").Value = "  -2.81%  "

# Row 49 - was PaxDollar, now InjectiveProtocol
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").Value = "'12.04"
$ws.Range("E49").Value = "  +0.55%  "

# Row 50 - was InjectiveProtocol, now PaxDollar
$ws.Range("B50").Value = "PaxDollar"
$ws.Range("C50").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D50").Value = "'0.998"
$ws.Range("E50").Value = "  -0.12%  "

# Row 51 - BabyDogeCoin
$ws.Range("D51").Value = "0.0₆0122"
$ws.Range("E51").Value = "  +0.90%  "
